$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = 25
$ws1.Range("H2").Value = 3.72
$ws1.Range("L2").Value = 0.92

$ws1.Range("D3").Value = 26
$ws1.Range("H3").Value = 2.62
$ws1.Range("L3").Value = 0.97

$ws1.Range("H4").Value = 1.62
$ws1.Range("L4").Value = 1.14

$ws1.Range("D5").Value = 28
$ws1.Range("H5").Value = 0.57
$ws1.Range("I5").Value = "Low"
$ws1.Range("L5").Value = 0.86

$ws1.Range("D6").Value = 26
$ws1.Range("L6").Value = 0.99

$ws1.Range("D7").Value = 24
$ws1.Range("L7").Value = 1.16

$ws1.Range("D8").Value = 23
$ws1.Range("L8").Value = 1.17

$ws1.Range("D9").Value = 23
$ws1.Range("L9").Value = 1.08

$ws1.Range("D10").Value = 24
$ws1.Range("L10").Value = 0.97

$ws1.Range("D11").Value = 22
$ws1.Range("L11").Value = 1.06

$ws1.Range("D12").Value = 21
$ws1.Range("L12").Value = 0.9399999999999999

$ws1.Range("D13").Value = 20
$ws1.Range("L13").Value = 0.85

$ws1.Range("D14").Value = 20
$ws1.Range("L14").Value = 1.16

$ws1.Range("D15").Value = 20
$ws1.Range("L15").Value = 1.15

$ws1.Range("D16").Value = 21
$ws1.Range("L16").Value = 0.97

$ws1.Range("D17").Value = 20
$ws1.Range("L17").Value = 1.19

# --- Sheet: Summary ---
# These cells hold numeric-looking text (e.g. "369") as strings in the
# original workbook. Force text format first so Excel doesn't silently
# convert the assigned string into a numeric value.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "369"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "201"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "105"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "28"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "20"
